$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '39.661.66'
$ws.Cells.Item(2, 5).Value = '  -1.14%  '
$ws.Cells.Item(3, 4).Value = '2.186.91'
$ws.Cells.Item(4, 5).Value = '  -0.10%  '
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = '290.12'
$ws.Cells.Item(5, 5).Value = '  -1.06%  '
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = '86.00'
$ws.Cells.Item(6, 5).Value = '  -1.65%  '
$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = '0.507'
$ws.Cells.Item(7, 5).Value = '  -1.84%  '
$ws.Cells.Item(8, 5).Value = '  -0.02%  '
$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = '0.461'
$ws.Cells.Item(9, 5).Value = '  -2.68%  '
$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = '29.96'
$ws.Cells.Item(10, 5).Value = '  -4.61%  '
$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = '50.01'
$ws.Cells.Item(11, 5).Value = '  +6.37%  '
$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = '0.0774'
$ws.Cells.Item(12, 5).Value = '  -2.28%  '
$ws.Cells.Item(13, 5).Value = '  +2.44%  '
$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = '6.37'
$ws.Cells.Item(14, 5).Value = '  -0.89%  '
$ws.Cells.Item(15, 4).Value = '2.526.63'
$ws.Cells.Item(15, 5).Value = '  -2.43%  '
$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = '13.59'
$ws.Cells.Item(16, 5).Value = '  -3.82%  '
$ws.Cells.Item(17, 4).Value = '2.146.57'
$ws.Cells.Item(17, 5).Value = '  -4.12%  '
$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = '0.721'
$ws.Cells.Item(18, 5).Value = '  -1.77%  '
$ws.Cells.Item(19, 4).Value = '39.576.56'
$ws.Cells.Item(19, 5).Value = '  -1.24%  '
$ws.Cells.Item(20, 4).Value = '0.0₃0876'
$ws.Cells.Item(20, 5).Value = '  -1.56%  '
$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = '11.09'
$ws.Cells.Item(21, 5).Value = '  -1.52%  '
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = '5.67'
$ws.Cells.Item(22, 5).Value = '  -3.10%  '
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = '64.91'
$ws.Cells.Item(23, 5).Value = '  -1.09%  '
$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = '235.36'
$ws.Cells.Item(24, 5).Value = '  -0.46%  '
$ws.Cells.Item(25, 5).Value = '  +0.14%  '
$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = '2.42'
$ws.Cells.Item(26, 5).Value = '  -2.46%  '
$ws.Cells.Item(27, 5).Value = '  -4.06%  '
$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = '22.86'
$ws.Cells.Item(28, 5).Value = '  -0.44%  '
$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = '2.14'
$ws.Cells.Item(29, 5).Value = '  -3.66%  '
$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = '9.10'
$ws.Cells.Item(30, 5).Value = '  -2.60%  '
$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = '155.75'
$ws.Cells.Item(31, 5).Value = '  +2.56%  '
$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = '30.99'
$ws.Cells.Item(32, 5).Value = '  -7.42%  '
$ws.Cells.Item(33, 5).Value = '  -0.16%  '
$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = '4.87'
$ws.Cells.Item(34, 5).Value = '  -1.58%  '
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = '0.0701'
$ws.Cells.Item(35, 5).Value = '  -2.99%  '
$ws.Cells.Item(36, 5).Value = '  -2.40%  '
$ws.Cells.Item(37, 5).Value = '  -0.94%  '
$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = '0.111'
$ws.Cells.Item(38, 5).Value = '  -0.59%  '
$ws.Cells.Item(39, 5).Value = '  -3.45%  '
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = '1.65'
$ws.Cells.Item(40, 5).Value = '  -4.49%  '
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = '14.91'
$ws.Cells.Item(42, 4).Value = '2.104.65'
$ws.Cells.Item(42, 5).Value = '  +1.30%  '
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = '3.69'
$ws.Cells.Item(43, 5).Value = '  -4.07%  '
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = '2.09'
$ws.Cells.Item(44, 5).Value = '  -1.09%  '
$ws.Cells.Item(45, 5).Value = '  -1.99%  '
$ws.Cells.Item(46, 5).Value = '  -2.42%  '
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = '17.06'
$ws.Cells.Item(47, 5).Value = '  -7.58%  '
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = '2.64'
$ws.Cells.Item(48, 5).Value = '  +0.91%  '
$ws.Cells.Item(49, 4).Value = '2.401.69'
$ws.Cells.Item(49, 5).Value = '  -1.78%  '
$ws.Cells.Item(50, 5).Value = '  -0.60%  '
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = '87.38'
$ws.Cells.Item(51, 5).Value = '  -2.41%  '
